$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the data that used to be in row 3 (HARRY), plus new values
# Force text formatting so "1011" stays a text value (like "0001" was), not a number
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1011"
$ws.Range("B2").Value = "HARRY"
$ws.Range("C2").Value = "19:57:23"
$ws.Range("D2").Value = "19:57:38"
$ws.Range("E2").Value = "0:00:15"
$ws.Range("F2").Value = "Present"
$ws.Range("G2").Value = "dataset/1011/1.jpg"

# Delete row 3 entirely (shifts rows up, removing the old row 3 data)
$ws.Rows(3).Delete()
